# Auto-generated script to update Fenrir_Profits market data values
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 43407116
$ws.Range("I98").Value = 46300784
$ws.Range("K98").Value = 46300784
$ws.Range("M98").Value = -46299286
$ws.Range("H116").Value = 275187.2
$ws.Range("I116").Value = 3931.1538
$ws.Range("J116").Value = 422117.53
$ws.Range("K116").Value = 3931.1538
$ws.Range("L116").Value = 422117.53
$ws.Range("M116").Value = -489.1538
$ws.Range("N116").Value = -429001.53
$ws.Range("H122").Value = 43407116
$ws.Range("I122").Value = 46300784
$ws.Range("K122").Value = 138902352
$ws.Range("M122").Value = -138899902
$ws.Range("H129").Value = 760.8571
$ws.Range("J129").Value = 1009.8788
$ws.Range("L129").Value = 3029.6364
$ws.Range("N129").Value = -13029.6364
$ws.Range("H132").Value = 20620624
$ws.Range("I132").Value = 23662850
$ws.Range("J132").Value = 1087.3334
$ws.Range("K132").Value = 70988550
$ws.Range("L132").Value = 3262.0002
$ws.Range("M132").Value = -70986020
$ws.Range("N132").Value = -8322.0002
$ws.Range("H135").Value = 3898.1086
$ws.Range("I135").Value = 4397.727
$ws.Range("J135").Value = 2629.8462
$ws.Range("K135").Value = 39579.543
$ws.Range("L135").Value = 23668.6158
$ws.Range("M135").Value = -37044.543
$ws.Range("N135").Value = -28738.6158
$ws.Range("H137").Value = 37479224
$ws.Range("I137").Value = 1191542.9
$ws.Range("J137").Value = 58824920
$ws.Range("K137").Value = 3574628.7
$ws.Range("L137").Value = 176474760
$ws.Range("M137").Value = -3572078.7
$ws.Range("N137").Value = -176479860
$ws.Range("H138").Value = 1298.15
$ws.Range("I138").Value = 668.8276
$ws.Range("J138").Value = 2167.2144
$ws.Range("K138").Value = 2006.4828
$ws.Range("L138").Value = 6501.6432
$ws.Range("M138").Value = 3133.5172
$ws.Range("N138").Value = -16781.6432
$ws.Range("H141").Value = 2414.7922
$ws.Range("I141").Value = 2181.324
$ws.Range("J141").Value = 5177.5
$ws.Range("K141").Value = 6543.972
$ws.Range("L141").Value = 15532.5
$ws.Range("M141").Value = -1363.972
$ws.Range("N141").Value = -25892.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2127.75
$ws.Range("I2").Value = 1837
$ws.Range("J2").Value = 3000
$ws.Range("K2").Value = 1837
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = -1724
$ws.Range("N2").Value = -3226
$ws.Range("H61").Value = 8133284.5
$ws.Range("I61").Value = 9262812
$ws.Range("J61").Value = 691.2
$ws.Range("K61").Value = 9262812
$ws.Range("L61").Value = 691.2
$ws.Range("M61").Value = -9262600
$ws.Range("N61").Value = -1115.2
$ws.Range("H74").Value = 639.65
$ws.Range("I74").Value = 445.93332
$ws.Range("J74").Value = 1220.8
$ws.Range("K74").Value = 445.93332
$ws.Range("L74").Value = 1220.8
$ws.Range("M74").Value = 428.06668
$ws.Range("N74").Value = -2968.8
$ws.Range("H77").Value = 639.65
$ws.Range("I77").Value = 445.93332
$ws.Range("J77").Value = 1220.8
$ws.Range("K77").Value = 2229.6666
$ws.Range("L77").Value = 6104
$ws.Range("M77").Value = 2138.3334
$ws.Range("N77").Value = -14840
$ws.Range("H116").Value = 2127.75
$ws.Range("I116").Value = 1837
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 1837
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = 457
$ws.Range("N116").Value = -7588
$ws.Range("H136").Value = 8133284.5
$ws.Range("I136").Value = 9262812
$ws.Range("J136").Value = 691.2
$ws.Range("K136").Value = 27788436
$ws.Range("L136").Value = 2073.6
$ws.Range("M136").Value = -27785886
$ws.Range("N136").Value = -7173.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2127.75
$ws.Range("I3").Value = 1837
$ws.Range("J3").Value = 3000
$ws.Range("K3").Value = 1837
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = -1723
$ws.Range("N3").Value = -3228
$ws.Range("H134").Value = 13344014
$ws.Range("I134").Value = 17557608
$ws.Range("J134").Value = 963.1667
$ws.Range("K134").Value = 52672824
$ws.Range("L134").Value = 2889.5001
$ws.Range("M134").Value = -52670289
$ws.Range("N134").Value = -7959.5001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6296206
$ws.Range("I31").Value = 7576816
$ws.Range("J31").Value = 35446
$ws.Range("K31").Value = 7576816
$ws.Range("L31").Value = 35446
$ws.Range("M31").Value = -7576521
$ws.Range("N31").Value = -36036
$ws.Range("H34").Value = 6296206
$ws.Range("I34").Value = 7576816
$ws.Range("J34").Value = 35446
$ws.Range("K34").Value = 7576816
$ws.Range("L34").Value = 35446
$ws.Range("M34").Value = -7576614
$ws.Range("N34").Value = -35850
$ws.Range("H132").Value = 8776333
$ws.Range("I132").Value = 10754457
$ws.Range("J132").Value = 16072
$ws.Range("K132").Value = 32263371
$ws.Range("L132").Value = 48216
$ws.Range("M132").Value = -32260841
$ws.Range("N132").Value = -53276
$ws.Range("H134").Value = 19590680
$ws.Range("I134").Value = 24510672
$ws.Range("J134").Value = 3908202
$ws.Range("K134").Value = 73532016
$ws.Range("L134").Value = 11724606
$ws.Range("M134").Value = -73529481
$ws.Range("N134").Value = -11729676

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 55558696
$ws.Range("I132").Value = 65973430
$ws.Range("J132").Value = 13423.667
$ws.Range("K132").Value = 197920290
$ws.Range("L132").Value = 40271.001
$ws.Range("M132").Value = -197917760
$ws.Range("N132").Value = -45331.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2899606.2
$ws.Range("I132").Value = 4546121
$ws.Range("J132").Value = 1740.92
$ws.Range("K132").Value = 13638363
$ws.Range("L132").Value = 5222.76
$ws.Range("M132").Value = -13635833
$ws.Range("N132").Value = -10282.76
$ws.Range("H136").Value = 15876063
$ws.Range("I136").Value = 22225602
$ws.Range("J136").Value = 2214.6667
$ws.Range("K136").Value = 66676806
$ws.Range("L136").Value = 6644.000100000001
$ws.Range("M136").Value = -66674256
$ws.Range("N136").Value = -11744.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 63204.5
$ws.Range("J46").Value = 63204.5
$ws.Range("L46").Value = 63204.5
$ws.Range("N46").Value = -63666.5
$ws.Range("H107").Value = 22727668
$ws.Range("I107").Value = 436.9375
$ws.Range("K107").Value = 1310.8125
$ws.Range("M107").Value = 609.1875
$ws.Range("H132").Value = 28568398
$ws.Range("I132").Value = 46155170
$ws.Range("J132").Value = 15119691
$ws.Range("K132").Value = 138465510
$ws.Range("L132").Value = 45359073
$ws.Range("M132").Value = -138462980
$ws.Range("N132").Value = -45364133
$ws.Range("H134").Value = 63204.5
$ws.Range("J134").Value = 63204.5
$ws.Range("L134").Value = 189613.5
$ws.Range("N134").Value = -194683.5
$ws.Range("H136").Value = 20408102
$ws.Range("I136").Value = 17856828
$ws.Range("J136").Value = 33334556
$ws.Range("K136").Value = 53570484
$ws.Range("L136").Value = 100003668
$ws.Range("M136").Value = -53567934
$ws.Range("N136").Value = -100008768
